# Add three more "in progress" rows to the Progress sheet, mirroring the
# existing rows for the word "ਲੇਪਨ" but recorded at a later "selected_at"
# timestamp (45914.24793988426).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Progress")

$dateFormat = $ws.Range("G2").NumberFormat

$word = "ਲੇਪਨ"
$wordKeyNorm = "ਲੇਪਨ"
$selectedAt = 45914.24793988426

$verses = @(
    "ਚੰਦਨ ਅਗਰ ਕਪੂਰ ਲੇਪਨ ਤਿਸੁ ਸੰਗੇ ਨਹੀ ਪ੍ਰੀਤਿ ॥",
    "ਜਟਾ ਭਸਮ ਲੇਪਨ ਕੀਆ ਕਹਾ ਗੁਫਾ ਮਹਿ ਬਾਸੁ ॥",
    "ਬਾਹਰਿ ਭਸਮ ਲੇਪਨ ਕਰੇ ਅੰਤਰਿ ਗੁਬਾਰੀ ॥"
)
$pageNumbers = @(1018, 1103, 1243)

for ($i = 0; $i -lt $verses.Length; $i++) {
    $row = 5 + $i

    $ws.Cells.Item($row, 1).Value = $word
    $ws.Cells.Item($row, 2).Value = $wordKeyNorm
    $ws.Cells.Item($row, 4).Value = $verses[$i]
    $ws.Cells.Item($row, 5).Value = $pageNumbers[$i]
    $ws.Cells.Item($row, 6).Value = $true
    $ws.Cells.Item($row, 7).Value = $selectedAt
    $ws.Cells.Item($row, 7).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 8).Value = "not started"
}
